# Insert a new paragraph "School will resume from Shawan 7th." right after the
# "Write a HTML code to create a hotel menu." paragraph (and before the
# trailing empty paragraph), matching the formatting shown in the target
# diff (centered, bold/italic/underlined Times New Roman 16pt, with "th"
# as a superscript run).

$d = $word.ActiveDocument

# Locate the "hotel menu" paragraph by its text content so the script does
# not depend on a hard-coded paragraph index.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Write a HTML code to create a hotel menu.*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'hotel menu' paragraph"
}

$hotelMenuPara = $d.Paragraphs.Item($targetIndex)

# Create a new, blank paragraph right after it.
[void]$hotelMenuPara.Range.InsertParagraphAfter()

# Re-fetch the freshly inserted (now-empty) paragraph.
$newPara = $d.Paragraphs.Item($targetIndex + 1)

# Populate it (and all of its formatting, including the superscript "th"
# run) in one shot via WordOpenXML/InsertXML so the exact run split and
# run/paragraph-mark formatting from the target revision is reproduced.
$xmlSnippet = '<?xml version="1.0"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>School will resume from Shawan 7</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$newPara.Range.InsertXML($xmlSnippet)
